$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Correction epaisseur inserts: Epaisseur_rotule (B12) 9 -> 6.75
$ws.Range("B12").Value = 6.75

# Move the active selection down to B13, as the author left it when saving
$ws.Range("B13").Select()
